{"js": "// Replace the two `source(file.path(root_dir, \"...\"))` /\n// `read_sav(file.path(root_dir, \"...\"))` calls with `find_file_up(...)`:\n//   file.path            -> find_file_up   (FunctionTok run)\n//   (root_dir, \" ... \")  -> (\" ... \")      (\"(root_dir, \" -> \"(\" in the NormalTok run)\nconst fnResults = context.document.body.search(\"file.path\", { matchCase: true });\nfnResults.load(\"text\");\nawait context.sync();\n\nfor (const r of fnResults.items) {\n  r.insertText(\"find_file_up\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst argResults = context.document.body.search(\"(root_dir, \", { matchCase: true });\nargResults.load(\"text\");\nawait context.sync();\n\nfor (const r of argResults.items) {\n  r.insertText(\"(\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Replace the two write-up paragraphs (Problem 13 / Bonferroni contrast, and\n// Problem 25 / Northeast vs Midwest & West contrast) with the new APA-style text.\nconst oldP14 =\n  \"Using the Bonferroni-adjusted threshold (\\u03b1 = .05/4), the city vs non-city contrast is significant in the Midwest and Northeast, but not significant in the South or West. Directions, effect sizes (d), and Bonferroni-adjusted p-values are shown above.\";\nconst newP14 =\n  \"APA-style write-up (Bonferroni-adjusted): In the Midwest, city dwellers scored higher than non-city dwellers, t(312) = 2.96, p_bonf = .013, d = 0.66. In the Northeast, non-city dwellers scored higher than city dwellers, t(312) = 7.23, p_bonf < .001, d = 1.62. The South and West comparisons were not significant (both p_bonf = 1.00).\";\n\nconst oldP25 =\n  \"Write-up: The Northeast vs Midwest & West contrast is significant for both city and non-city dwellers, and the interaction contrast is also significant. Report the direction of effects, t ratios, p-values, and Cohen\\u2019s d from the contrast table above, using the Scheffe and Bonferroni critical values for reference.\";\nconst newP25 =\n  \"APA-style write-up: For city dwellers, the Northeast vs (Midwest & West) contrast was significant, t(312) = 5.61, p < .001, d = 1.09, indicating lower Northeast scores than the average of Midwest and West for city residents (contrast estimate = -1.03). For non-city dwellers, the same contrast was significant, t(312) = 3.83, p = .0002, d = 0.74, indicating higher Northeast scores than the average of Midwest and West (estimate = 0.70). The interaction contrast was significant, t(312) = 6.67, p < .001, d = 1.83, showing that the Northeast advantage differs sharply by city status.\";\n\nconst p14Results = context.document.body.search(oldP14, { matchCase: true });\np14Results.load(\"text\");\nawait context.sync();\nfor (const r of p14Results.items) {\n  r.insertText(newP14, Word.InsertLocation.replace);\n}\n\nconst p25Results = context.document.body.search(oldP25, { matchCase: true });\np25Results.load(\"text\");\nawait context.sync();\nfor (const r of p25Results.items) {\n  r.insertText(newP25, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Replace the two `file.path` function-call sites with `find_file_up`,\n#    and drop the `root_dir, ` leading argument text (keeping the opening\n#    paren) so `file.path(root_dir, \"...\")` becomes `find_file_up(\"...\")`.\n$find = $d.Content.Find\n$find.Text = \"file.path\"\n$find.Replacement.Text = \"find_file_up\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \"(root_dir, \"\n$find2.Replacement.Text = \"(\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) Replace the Problem 13/14 write-up paragraph with the new APA-style text.\n$oldP14 = \"Using the Bonferroni-adjusted threshold (\" + [char]0x03B1 + \" = .05/4), the city vs non-city contrast is significant in the Midwest and Northeast, but not significant in the South or West. Directions, effect sizes (d), and Bonferroni-adjusted p-values are shown above.\"\n$newP14 = \"APA-style write-up (Bonferroni-adjusted): In the Midwest, city dwellers scored higher than non-city dwellers, t(312) = 2.96, p_bonf = .013, d = 0.66. In the Northeast, non-city dwellers scored higher than city dwellers, t(312) = 7.23, p_bonf < .001, d = 1.62. The South and West comparisons were not significant (both p_bonf = 1.00).\"\n\n$find3 = $d.Content.Find\n$find3.Text = $oldP14\n$find3.Replacement.Text = $newP14\n$find3.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 3) Replace the Problem 25 write-up paragraph with the new APA-style text.\n$oldP25 = \"Write-up: The Northeast vs Midwest & West contrast is significant for both city and non-city dwellers, and the interaction contrast is also significant. Report the direction of effects, t ratios, p-values, and Cohen\" + [char]0x2019 + \"s d from the contrast table above, using the Scheffe and Bonferroni critical values for reference.\"\n$newP25 = \"APA-style write-up: For city dwellers, the Northeast vs (Midwest & West) contrast was significant, t(312) = 5.61, p < .001, d = 1.09, indicating lower Northeast scores than the average of Midwest and West for city residents (contrast estimate = -1.03). For non-city dwellers, the same contrast was significant, t(312) = 3.83, p = .0002, d = 0.74, indicating higher Northeast scores than the average of Midwest and West (estimate = 0.70). The interaction contrast was significant, t(312) = 6.67, p < .001, d = 1.83, showing that the Northeast advantage differs sharply by city status.\"\n\n$find4 = $d.Content.Find\n$find4.Text = $oldP25\n$find4.Replacement.Text = $newP25\n$find4.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
